$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 566; existing rows 566:640 shift down to 567:641
$ws.Rows("566:566").Insert()

# Populate the newly inserted row 566 with the new weekly price observation
$ws.Cells.Item(566, 1).Value2 = 9
$ws.Cells.Item(566, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(566, 3).Value2 = "Metropolitana"
$ws.Cells.Item(566, 4).Value2 = 45127
$ws.Cells.Item(566, 5).Value2 = 13
$ws.Cells.Item(566, 6).Value2 = 100112012
$ws.Cells.Item(566, 7).Value2 = "Espinaca"
$ws.Cells.Item(566, 8).Value2 = "Sin especificar"
$ws.Cells.Item(566, 9).Value2 = "Primera"
$ws.Cells.Item(566, 10).Value2 = 160
$ws.Cells.Item(566, 11).Value2 = 7000
$ws.Cells.Item(566, 12).Value2 = 8000
$ws.Cells.Item(566, 13).Value2 = 7500
$ws.Cells.Item(566, 14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(566, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(566, 16).Value2 = 750
$ws.Cells.Item(566, 17).Value2 = 10
$ws.Cells.Item(566, 18).Value2 = "Hortaliza"
